$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new summary row (18) that mirrors row 17's layout, with fresh
# "Q. 0.0" totals in the last four columns.
$ws.Range("A18").Value = "21 de octubre del 2024"
$ws.Range("B18").Value = "Resumen"
$ws.Range("C18").Value = "total"
$ws.Range("D18").Value = "del"
$ws.Range("E18").Value = "dia"
$ws.Range("F18").Value = "-"
$ws.Range("G18").Value = "Q. 0.0"
$ws.Range("H18").Value = "Q. 0.0"
$ws.Range("I18").Value = "Q. 0.0"
$ws.Range("J18").Value = "Q. 0.0"
